$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "done, using a global messenger system"
Write-Host "t1"
